$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.546.55'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.51%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.220.09'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.91%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.77'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.17'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.98%  '

$ws.Range('E7').Value = '  -3.85%  '

$ws.Range('E8').Value = '  +0.15%  '

$ws.Range('E9').Value = '  -7.48%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.62'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -8.37%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0817'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.75%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.32'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -7.18%  '

$ws.Range('E13').Value = '  -3.31%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.558.52'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.71%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.833'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -5.37%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.214.78'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.24%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.02'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.08%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.410.49'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.58%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.90'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -10.90%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0957'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.96%  '

$ws.Range('E21').Value = '  -6.33%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.04'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.83%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.28'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.55%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.96'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -8.35%  '

$ws.Range('E25').Value = '  -8.37%  '

$ws.Range('E26').Value = '  +0.18%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.91'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.56%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.77%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.34'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -6.55%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '159.14'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.09%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.90'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -9.95%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.75'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.86%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0822'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -7.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.65'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.39%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.11'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.07%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.86'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -10.58%  '

$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.107'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.58%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.117'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.72%  '

$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.28'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.72%  '

$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.52'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -9.74%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.98'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -12.90%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0304'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -7.18%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.12%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.695.36'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.93%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '83.53'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.35%  '

$ws.Range('E46').Value = '  -7.68%  '

$ws.Range('E47').Value = '  -6.55%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.11'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.89%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.62'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.59%  '

$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.47'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -6.34%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.77'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.75%  '
